$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 88.60676470588236

$ws.Range("I39:I71").Value = $newValue
